$d = $word.ActiveDocument

# Target paragraph sequence describing the full document body after the
# edit: a Title, three Heading1 sections (QUALITY / GENRE / DIVERSITY OF
# LOCATION / PROFIT) each followed by body-text paragraphs, separated by
# blank spacer paragraphs, exactly mirroring the authored writeup.
$paras = @(
    @{ style = "Normal"; text = "" },
    @{ style = "Title"; text = "Findings Report" },
    @{ style = "Normal"; text = "" },
    @{ style = "Heading1"; text = "QUALITY" },
    @{ style = "Normal"; text = "Our team has retrieved, a total of 5325 movies from 2000-2023 that are available in Netflix, Stan and Amazon Prime Video. In terms of variety, Amazon and Netflix provider viewers in Australia with far more options, with a large focus on recent releases. Across the three platforms, 2223 movies we’re release between 2015 to 2019." },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "By comparing this volume with the IMDb score rating for each film we can start to notice trends. While Amazon has saturated the market with options, it doesn’t necessarily equate to quality programming. Netflix and Stan both trend higher when it comes to overall voting scores, with Stan averaging higher overall for the 20-year period." },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "Netflix boasts 646 more available title to watch than Stan, while only having only a -0.2 rating difference in terms of quality (Stan 6.5, Netflix 6.3)" },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "In this aspect, it’s fair to argue that Netflix provides the best options for consumers wanting a large amount of quality viewing. It may be a strategy of streaming platforms to target audiences with preference `"trending`" movies within the decade. Additionally, there is a significant decrease of available movies from 2020 onwards and it may be due to the pandemic restriction on movie production." },
    @{ style = "Normal"; text = "" },
    @{ style = "Heading1"; text = "GENRE" },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "Using boxplots, we visualise the genre distribution of the top 100 rated movies available in the streaming platforms. Interestingly, in grouped year 2015-19, Drama and Comedy movies have the most outliers from both ends of IMDB ratings. The plots per year group shows several outliers of rating score from different genres and this may be due to the subjective preference and number of voters from the IMDB data file." },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "Using barplot, the figure shows comparison of movie genre across three different platforms: Netflix, Stan and Amazon Prime Video. This easily showcases the quality of available film categories which may be helpful for audiences to decide what streaming platform to use depending on their movie preference. For example, if Person A only prefers science-fiction movies, Stan has the most high rated movies in that genre." },
    @{ style = "Normal"; text = "" },
    @{ style = "Heading1"; text = "DIVERSITY OF LOCATION" },
    @{ style = "Normal"; text = "By plotting the location of the films country of origin, we can determine how diverse the streaming platforms library is." },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "(Size of the circle determines the volume of films from that location, where intensity of color shows average ratings.)" },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "Using the visuals, we can see all three providers are US dominated. Amazon and Netflix boast a far broader range of countries than Stan does, but Stan makes up for this by having the higher quality films on average per country available. " },
    @{ style = "Normal"; text = "" },
    @{ style = "Heading1"; text = "PROFIT" },
    @{ style = "Normal"; text = "The graph for the Top 500 highest budgeted films with their IMDB ratings shows the genre distribution of movies. It is interesting to note that majority of genres have a budget less than 1 million in USD with rating score of 5 to 8. The movie genres in the highest end of budget spectrum with high ratings are Fantasy, Sci-fi and Adventure film categories (e.g. Justice League, Spiderman, The Hobbit). " },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "We can also map out the ‘profit’ of the films available (not to be confused with the profit of the streaming platforms) to identify how many “successful” films each streaming service hosts. " },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "Split over over the 20 year period, we can see a sharp decline in overall film profits starting at around 2015. This can be directly correlated with the introduction of streaming services, highlighting that larger box office productions started taking a hit when consumers were given the option to watch from home. This is further proven by Netflix recording a higher profit from 2015 - 2020, as Netflix dominated the market with their ‘Original Content’" },
    @{ style = "Normal"; text = "" },
    @{ style = "Normal"; text = "It may also highlight the change in peoples viewing preference from Films to TV Series." }
)

# The document starts with 3 paragraphs (Heading1 title / blank / Overview).
# Keep paragraph #1 and the final paragraph (which must keep carrying the
# sectPr / final paragraph mark) in place, and grow the body to the full
# target length by repeatedly inserting new blank paragraphs right after
# paragraph #1 -- this naturally pushes the original #2/#3 paragraphs down
# while the very last paragraph of the document stays last.
$target = $paras.Length
$idx = 1
while ($d.Paragraphs.Count -lt $target) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
}

# Every target paragraph slot now exists; set style + text for each by index.
for ($i = 1; $i -le $paras.Length; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Style = $paras[$i - 1].style
    $p.Range.Text = $paras[$i - 1].text
}

Write-Host "Paragraph count: $($d.Paragraphs.Count)"
